$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "new speed" column: header in E1, single data value in E3
$ws.Range("E1").Value = "new speed"
$ws.Range("E3").Value = 9

# Leave the selection where the author left it after entering the data
$ws.Range("E4").Select()
